$wb = $excel.ActiveWorkbook

# --- Sheet "UCM" (sheet1) ---
$ucm = $wb.Worksheets.Item("UCM")
$ucm.Range("H8").Value = 20
$ucm.Range("H9").Select()

# --- Sheet "R45" (sheet2) ---
$r45 = $wb.Worksheets.Item("R45")

# Fill in the first fold table (rows 6-7)
$r45.Range("B6").Value = 143.42608537000999
$r45.Range("C6").Value = 0.088832819927011003
$r45.Range("D6").Value = 0.97174603174603102

$r45.Range("B7").Value = 143.85543215891801
$r45.Range("C7").Value = 0.078710013302752194
$r45.Range("D7").Value = 0.97593650793650699

# Fill in the second fold table (rows 13-14)
$r45.Range("B13").Value = 218.71668514295001
$r45.Range("C13").Value = 0.041701368287886398
$r45.Range("D13").Value = 0.98869841269841197

$r45.Range("B14").Value = 217.143963020993
$r45.Range("C14").Value = 0.032198057287765001
$r45.Range("D14").Value = 0.99060317460317404

$r45.Range("G13:J18").Select()

$wb.Save()
